# Append a new data row (row 49) to the worksheet, duplicating the last
# existing row (row 48), which is what the Adafruit IO export produced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 48
$newRow = $lastRow + 1

# Copy the last row and paste it into the new row so the new cells keep
# the same text representation / cell formatting as the source row
# (all values here are plain text, e.g. "25" rather than a number 25).
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial()
$excel.CutCopyMode = $false
